$d = $word.ActiveDocument
$p2 = $d.Paragraphs(2).Range
$newXml = @'
<w:p w:rsidR="0020009C" w:rsidRDefault="00835889"><w:r><w:tab/></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>theme</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>, design choices</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">Algorithm (abstract) (~2 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mins</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r><w:r><w:br/></w:r><w:r><w:tab/><w:t>why use this algorithm, why it's better/easier than others</w:t></w:r><w:r><w:br/><w:t>Algorithm (implementation) (~2mins)</w:t></w:r><w:r><w:br/></w:r><w:r><w:tab/><w:t>what changed to make the algorithm work</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">demonstration (~3 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mins</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r><w:r><w:br/></w:r><w:r><w:tab/></w:r><w:r><w:t xml:space="preserve">Sample input is going to be starting location:12, length 1000, and difficulty 10. Use a grid to show how the paths are working and what sort of paths we get on a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>google</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> map slide with all three different routes we got displayed. Explain how we go from nodes to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>latt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> long pairs.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>explain</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> plans on what to do with the difficulty.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:br/></w:r><w:r><w:br/><w:t xml:space="preserve">GOAL: 9 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mins</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$p2.InsertXML($newXml)
Write-Output $d.Paragraphs.Count
Write-Output $d.Paragraphs(2).Range.Text
